$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "Stocks"

# Duplicate the sheet (keeps identical sheetPr/sheetFormatPr/pageMargins/styles)
# so the new "Crypto" tab starts from the same base as "Stocks".
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Crypto"

# ---- Stocks sheet: overwrite rows 2-11 in place (keeps existing cell styles) ----
$ws1.Range("A2").Value = "AAPL"
$ws1.Range("B2").Value = [double]"0.001054895583622141"
$ws1.Range("C2").Value = [double]"0.01847867590763628"
$ws1.Range("A3").Value = "AMZN"
$ws1.Range("B3").Value = [double]"0.001178874713675017"
$ws1.Range("C3").Value = [double]"0.02102636305894617"
$ws1.Range("A4").Value = "BABA"
$ws1.Range("B4").Value = [double]"0.0002677426938746988"
$ws1.Range("C4").Value = [double]"0.02667360023865977"
$ws1.Range("A5").Value = "GOOG"
$ws1.Range("B5").Value = [double]"0.0008966390261401635"
$ws1.Range("C5").Value = [double]"0.0179986699792909"
$ws1.Range("A6").Value = "GOOGL"
$ws1.Range("B6").Value = [double]"0.0008872120707870734"
$ws1.Range("C6").Value = [double]"0.01795765358692025"
$ws1.Range("A7").Value = "JNJ"
$ws1.Range("B7").Value = [double]"0.0003566857385153873"
$ws1.Range("C7").Value = [double]"0.01162558218396622"
$ws1.Range("A8").Value = "JPM"
$ws1.Range("B8").Value = [double]"0.0006470654447110065"
$ws1.Range("C8").Value = [double]"0.01767014212288854"
$ws1.Range("A9").Value = "MSFT"
$ws1.Range("B9").Value = [double]"0.001088953473443014"
$ws1.Range("C9").Value = [double]"0.01764321156525108"
$ws1.Range("A10").Value = "V"
$ws1.Range("B10").Value = [double]"0.000720909043013748"
$ws1.Range("C10").Value = [double]"0.01599213059296209"
$ws1.Range("A11").Value = "VOD"
$ws1.Range("B11").Value = [double]"-0.0001898836618308103"
$ws1.Range("C11").Value = [double]"0.01724309977823213"
# Remove the remaining rows (former crypto/leftover stock rows)
$ws1.Range("A12:C26").EntireRow.Delete()

# ---- Crypto sheet: overwrite rows 2-16 in place (keeps existing cell styles) ----
$ws2.Range("A2").Value = "ADA-USD"
$ws2.Range("B2").Value = [double]"0.002109322907709398"
$ws2.Range("C2").Value = [double]"0.05385911692947309"
$ws2.Range("A3").Value = "BNB-USD"
$ws2.Range("B3").Value = [double]"0.003300274737707548"
$ws2.Range("C3").Value = [double]"0.05275003261595376"
$ws2.Range("A4").Value = "BTC-USD"
$ws2.Range("B4").Value = [double]"0.001338890310131067"
$ws2.Range("C4").Value = [double]"0.03414030084679871"
$ws2.Range("A5").Value = "DAI-USD"
$ws2.Range("B5").Value = [double]"-2.317790536334618e-05"
$ws2.Range("C5").Value = [double]"0.002147598150660589"
$ws2.Range("A6").Value = "DOGE-USD"
$ws2.Range("B6").Value = [double]"0.007118002616628462"
$ws2.Range("C6").Value = [double]"0.1316804669561277"
$ws2.Range("A7").Value = "ETH-USD"
$ws2.Range("B7").Value = [double]"0.002288183890886724"
$ws2.Range("C7").Value = [double]"0.0450493043897679"
$ws2.Range("A8").Value = "LINK-USD"
$ws2.Range("B8").Value = [double]"0.001057863688019033"
$ws2.Range("C8").Value = [double]"0.05777973265314883"
$ws2.Range("A9").Value = "MATIC-USD"
$ws2.Range("B9").Value = [double]"0.005231289252545505"
$ws2.Range("C9").Value = [double]"0.07498958764790675"
$ws2.Range("A10").Value = "SOL-USD"
$ws2.Range("B10").Value = [double]"0.004044725884465609"
$ws2.Range("C10").Value = [double]"0.07182484297952729"
$ws2.Range("A11").Value = "TON-USD"
$ws2.Range("B11").Value = [double]"0.0006277791916024"
$ws2.Range("C11").Value = [double]"0.06824254427809426"
$ws2.Range("A12").Value = "TRX-USD"
$ws2.Range("B12").Value = [double]"0.002231566251630027"
$ws2.Range("C12").Value = [double]"0.04737623518475711"
$ws2.Range("A13").Value = "USDC-USD"
$ws2.Range("B13").Value = [double]"-1.439140654660258e-06"
$ws2.Range("C13").Value = [double]"0.001241952029302097"
$ws2.Range("A14").Value = "USDT-USD"
$ws2.Range("B14").Value = [double]"-2.050810642745787e-06"
$ws2.Range("C14").Value = [double]"0.0007281035275920827"
$ws2.Range("A15").Value = "WBTC-USD"
$ws2.Range("B15").Value = [double]"0.001281691832933432"
$ws2.Range("C15").Value = [double]"0.03417820766614947"
$ws2.Range("A16").Value = "XRP-USD"
$ws2.Range("B16").Value = [double]"0.002536448557910965"
$ws2.Range("C16").Value = [double]"0.06508533059648003"
# Remove the remaining rows (former leftover rows)
$ws2.Range("A17:C26").EntireRow.Delete()

$ws1.Select()
